$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "#"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Round 1"
$ws.Range("D1").Value = "Round 2"
$ws.Range("E1").Value = "Total"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Kai Ze Lim"
$ws.Range("C2").Value = 42900
$ws.Range("D2").Value = 57400
$ws.Range("E2").Value = 100300
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Long Ha"
$ws.Range("C3").Value = 2500
$ws.Range("D3").Value = 84700
$ws.Range("E3").Value = 87200
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Zeren Shen"
$ws.Range("C4").Value = 115200
$ws.Range("D4").Value = -30000
$ws.Range("E4").Value = 85200
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Sanjay Ravichandran "
$ws.Range("C5").Value = -30000
$ws.Range("D5").Value = 86000
$ws.Range("E5").Value = 56000
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Amery Caleb Atinon"
$ws.Range("C6").Value = 5000
$ws.Range("D6").Value = 41800
$ws.Range("E6").Value = 46800
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Yong Jun Loo"
$ws.Range("C7").Value = 67400
$ws.Range("D7").Value = -30000
$ws.Range("E7").Value = 37400
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "YuCheng Chien"
$ws.Range("C8").Value = -30000
$ws.Range("D8").Value = 63300
$ws.Range("E8").Value = 33300
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Tanmay Shewale "
$ws.Range("C9").Value = 54400
$ws.Range("D9").Value = -30000
$ws.Range("E9").Value = 24400
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Yunlin P‚ÄÜan"
$ws.Range("C10").Value = -30000
$ws.Range("D10").Value = 50000
$ws.Range("E10").Value = 20000
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Sean Bellato"
$ws.Range("C11").Value = 45700
$ws.Range("D11").Value = -30000
$ws.Range("E11").Value = 15700
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Eugene Yap"
$ws.Range("C12").Value = -30000
$ws.Range("D12").Value = 38400
$ws.Range("E12").Value = 8400
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Hong Sheng Quah"
$ws.Range("C13").Value = 33900
$ws.Range("D13").Value = -30000
$ws.Range("E13").Value = 3900
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Kevin Yu"
$ws.Range("C14").Value = -30000
$ws.Range("D14").Value = 30900
$ws.Range("E14").Value = 900
$ws.Range("A15").Value = 15
$ws.Range("B15").Value = "Zetong Zhang"
$ws.Range("C15").Value = 22100
$ws.Range("D15").Value = -30000
$ws.Range("E15").Value = -7900
$ws.Range("A16").Value = 16
$ws.Range("B16").Value = "Hao Duong"
$ws.Range("C16").Value = 18900
$ws.Range("D16").Value = -30000
$ws.Range("E16").Value = -11100
$ws.Range("A17").Value = 17
$ws.Range("B17").Value = "Ray chng"
$ws.Range("C17").Value = 18400
$ws.Range("D17").Value = -30000
$ws.Range("E17").Value = -11600
$ws.Range("A18").Value = 18
$ws.Range("B18").Value = "Cheng Sheng La "
$ws.Range("C18").Value = 10100
$ws.Range("D18").Value = -30000
$ws.Range("E18").Value = -19900
$ws.Range("A19").Value = 19
$ws.Range("B19").Value = "Justin Qiang"
$ws.Range("C19").Value = 8300
$ws.Range("D19").Value = -30000
$ws.Range("E19").Value = -21700
$ws.Range("A20").Value = 20
$ws.Range("B20").Value = "Ha Hwei Keat"
$ws.Range("C20").Value = 8200
$ws.Range("D20").Value = -30000
$ws.Range("E20").Value = -21800
$ws.Range("A21").Value = 21
$ws.Range("B21").Value = "Vik Ayyasamy Sivakumar"
$ws.Range("C21").Value = 7700
$ws.Range("D21").Value = -30000
$ws.Range("E21").Value = -22300
$ws.Range("A22").Value = 22
$ws.Range("B22").Value = "Owen Vandegraaff "
$ws.Range("C22").Value = 5600
$ws.Range("D22").Value = -30000
$ws.Range("E22").Value = -24400
$ws.Range("A23").Value = 23
$ws.Range("B23").Value = "Kodai Shichida"
$ws.Range("C23").Value = 2700
$ws.Range("D23").Value = -30000
$ws.Range("E23").Value = -27300
$ws.Range("A24").Value = 24
$ws.Range("B24").Value = "Sean Tan Yuheng "
$ws.Range("C24").Value = -30000
$ws.Range("D24").Value = -3100
$ws.Range("E24").Value = -33100
$ws.Range("A25").Value = 25
$ws.Range("B25").Value = "Josh Burke"
$ws.Range("C25").Value = -3600
$ws.Range("D25").Value = -30000
$ws.Range("E25").Value = -33600
$ws.Range("A26").Value = 26
$ws.Range("B26").Value = "Edwin Zou"
$ws.Range("C26").Value = -4700
$ws.Range("D26").Value = -30000
$ws.Range("E26").Value = -34700
$ws.Range("A27").Value = 27
$ws.Range("B27").Value = "Will Ho"
$ws.Range("C27").Value = -5200
$ws.Range("D27").Value = -30000
$ws.Range("E27").Value = -35200
$ws.Range("A28").Value = 28
$ws.Range("B28").Value = "Monil Bhatt"
$ws.Range("C28").Value = -6400
$ws.Range("D28").Value = -30000
$ws.Range("E28").Value = -36400
$ws.Range("A29").Value = 29
$ws.Range("B29").Value = "Justin Alexander T. Sy"
$ws.Range("C29").Value = -30000
$ws.Range("D29").Value = -18900
$ws.Range("E29").Value = -48900
$ws.Range("A30").Value = 30
$ws.Range("B30").Value = "Ashwin Seshadari "
$ws.Range("C30").Value = -20900
$ws.Range("D30").Value = -30000
$ws.Range("E30").Value = -50900
$ws.Range("A31").Value = 31
$ws.Range("B31").Value = "Kar Way Tan"
$ws.Range("C31").Value = -21600
$ws.Range("D31").Value = -30000
$ws.Range("E31").Value = -51600
$ws.Range("A32").Value = 32
$ws.Range("B32").Value = "Sinha Dayarathne"
$ws.Range("C32").Value = -30000
$ws.Range("D32").Value = -30000
$ws.Range("E32").Value = -60000
$ws.Range("A33").Value = 33
$ws.Range("B33").Value = "Deniz uragun"
$ws.Range("C33").Value = -30000
$ws.Range("D33").Value = -30000
$ws.Range("E33").Value = -60000
$ws.Range("A34").Value = 34
$ws.Range("B34").Value = "Goutham Peddireddy"
$ws.Range("C34").Value = -30000
$ws.Range("D34").Value = -30000
$ws.Range("E34").Value = -60000
$ws.Range("A35").Value = 35
$ws.Range("B35").Value = "Henry zhang"
$ws.Range("C35").Value = -30000
$ws.Range("D35").Value = -30000
$ws.Range("E35").Value = -60000
$ws.Range("A36").Value = 36
$ws.Range("B36").Value = "Edward Howells"
$ws.Range("C36").Value = -30000
$ws.Range("D36").Value = -30000
$ws.Range("E36").Value = -60000
$ws.Range("A37").Value = 37
$ws.Range("B37").Value = "Coco Yan"
$ws.Range("C37").Value = -30000
$ws.Range("D37").Value = -30000
$ws.Range("E37").Value = -60000
$ws.Range("A38").Value = 38
$ws.Range("B38").Value = "Ryan Cheng Hong Khoo"
$ws.Range("C38").Value = -30000
$ws.Range("D38").Value = -30000
$ws.Range("E38").Value = -60000
$ws.Range("A39").Value = 39
$ws.Range("B39").Value = "Lucy Liu "
$ws.Range("C39").Value = -30000
$ws.Range("D39").Value = -30000
$ws.Range("E39").Value = -60000
$ws.Range("A40").Value = 40
$ws.Range("B40").Value = "xinkai zhou"
$ws.Range("C40").Value = -30000
$ws.Range("D40").Value = -30000
$ws.Range("E40").Value = -60000
$ws.Range("A41").Value = 41
$ws.Range("B41").Value = "Yeo Zheng Xian"
$ws.Range("C41").Value = -30000
$ws.Range("D41").Value = -30000
$ws.Range("E41").Value = -60000
$ws.Range("A42").Value = 42
$ws.Range("B42").Value = "Qing Ye"
$ws.Range("C42").Value = -30000
$ws.Range("D42").Value = -30000
$ws.Range("E42").Value = -60000
$ws.Range("A43").Value = 43
$ws.Range("B43").Value = "Vinay Panicker"
$ws.Range("C43").Value = -30000
$ws.Range("D43").Value = -30000
$ws.Range("E43").Value = -60000
$ws.Range("A44").Value = 44
$ws.Range("B44").Value = "Ashwin Seshadari "
$ws.Range("C44").Value = -30000
$ws.Range("D44").Value = -30000
$ws.Range("E44").Value = -60000
$ws.Range("A45").Value = 45
$ws.Range("B45").Value = "James Patterson"
$ws.Range("C45").Value = -30000
$ws.Range("D45").Value = -30000
$ws.Range("E45").Value = -60000
$ws.Range("A46").Value = 46
$ws.Range("B46").Value = "Keji Yan"
$ws.Range("C46").Value = -30000
$ws.Range("D46").Value = -30000
$ws.Range("E46").Value = -60000
$ws.Range("A47").Value = 47
$ws.Range("B47").Value = "Nivethan Iyer"
$ws.Range("C47").Value = -30000
$ws.Range("D47").Value = -30000
$ws.Range("E47").Value = -60000

# Update selection to match target state (row 15 selected)
$ws.Range("A15:XFD15").Select()
